# "new results. 2021/06/07 23:50"
# Update two accuracy values in the SEED results sheet; the B32 average
# (=AVERAGE(B2:B31)) recalculates automatically from these.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B15").Value = 0.9212
$ws.Range("B26").Value = 0.922

# Scroll the view down a few rows (topLeftCell A11 -> A14), matching the
# author's saved scroll position.
$excel.ActiveWindow.ScrollRow = 14
